$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: update carrier and collapse detail columns into "mixed"
$ws.Range("C14").Value = "SAIA"
$ws.Range("F14").Value = "mixed"
$ws.Range("G14").Value = "mixed"
$ws.Range("H14").Value = "mixed"
$ws.Range("I14").Value = "mixed"

# L14/M14 need to become genuine numeric 1 (not text "1"), while keeping their
# existing style (s="9", which carries a Text number format). A direct
# .Value assignment on a Text-formatted cell is stored as text, so:
#  1) temporarily switch the number format to General
#  2) assign the numeric value
#  3) restore the original formatting by pasting formats from a cell that
#     still carries the original style (J14 keeps style 9 and is untouched)
$ws.Range("L14:M14").NumberFormat = "General"
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 1
$ws.Range("J14").Copy()
$ws.Range("L14:M14").PasteSpecial(-4122)

# Rows 15-22: clear all contents (cells become empty)
$ws.Range("A15:M22").ClearContents()
